$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.447.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.14%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.557.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.71%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'519.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'142.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.84%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.33%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.570.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.48%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.12%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.19%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -3.07%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.24%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.010.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.64%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'57.439.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.14%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'20.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.93%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.60%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.569.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.85%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'334.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.17%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.35%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.42%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.01%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'65.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.82%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.63%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.401"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -5.25%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.14%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.678.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'6.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.58%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -7.10%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -6.78%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.99%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.59%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'149.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.28%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.37%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.16%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.844"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -9.69%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'36.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.98%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E41").Value = "'  -1.08%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.25%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.13%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'268.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.70%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0956"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.09%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.13%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -4.02%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -4.32%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -3.21%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'RenderToken"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'4.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Maker"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'1.960.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.83%  "
$ws.Range("E51").Style = "Normal"
